$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 407) holds the "Förändrad" date serial value.
# All of these cells currently contain 45186 and must be updated to 45188.
$ws.Range("C2:C407").Value = 45188
